$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Semana 36 de 2025: add week 35 (col AL) and week 36 (col AM) ---
# Extend the weekly header row with the two new epidemiological week numbers,
# stored as text like the existing week-number headers (column D.."34" etc.).
$ws.Range("AL1:AM1").NumberFormat = "@"
$ws.Range("AL1").Value = "35"
$ws.Range("AM1").Value = "36"

# Copy the header style (bold, centered, General number format) from the previous week
# column (AK1) so the two new header cells look identical to the rest of row 1.
$ws.Range("AK1").Copy()
$ws.Range("AL1:AM1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New weekly report counts for week 35 (AL) and week 36 (AM), one value per UPGD row,
# plus a few corrected historical values (row 28, row 54) that came in with this update.
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = 0
$ws.Range("AL4").Value = 0
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = 0
$ws.Range("AL6").Value = 12
$ws.Range("AM6").Value = 24
$ws.Range("AL7").Value = 3
$ws.Range("AM7").Value = 1
$ws.Range("AL8").Value = 21
$ws.Range("AM8").Value = 15
$ws.Range("AL9").Value = 0
$ws.Range("AM9").Value = 0
$ws.Range("AL10").Value = 0
$ws.Range("AM10").Value = 0
$ws.Range("AM11").Value = 0
$ws.Range("AM12").Value = 0
$ws.Range("AL13").Value = 0
$ws.Range("AM13").Value = 0
$ws.Range("AL14").Value = 0
$ws.Range("AM14").Value = 0
$ws.Range("AM15").Value = 0
$ws.Range("AL16").Value = 0
$ws.Range("AM16").Value = 0
$ws.Range("AL17").Value = 0
$ws.Range("AM17").Value = 0
$ws.Range("AM18").Value = 0
$ws.Range("AM19").Value = 0
$ws.Range("AL22").Value = 0
$ws.Range("AM22").Value = 0
$ws.Range("AL23").Value = 0
$ws.Range("AM23").Value = 0
$ws.Range("AL25").Value = 1
$ws.Range("AM25").Value = 1
$ws.Range("AL26").Value = 0
$ws.Range("AM27").Value = 0
$ws.Range("AL28").Value = 5
$ws.Range("AM28").Value = 6
$ws.Range("O28").Value = 1
$ws.Range("Q28").Value = 1
$ws.Range("R28").Value = 2
$ws.Range("S28").Value = 2
$ws.Range("AL29").Value = 3
$ws.Range("AM29").Value = 3
$ws.Range("AL30").Value = 26
$ws.Range("AM30").Value = 17
$ws.Range("AL31").Value = 0
$ws.Range("AL35").Value = 16
$ws.Range("AL36").Value = 0
$ws.Range("AM36").Value = 0
$ws.Range("AL37").Value = 0
$ws.Range("AM37").Value = 0
$ws.Range("AL38").Value = 0
$ws.Range("AM38").Value = 0
$ws.Range("AL40").Value = 0
$ws.Range("AM40").Value = 0
$ws.Range("AL41").Value = 0
$ws.Range("AM41").Value = 0
$ws.Range("AL42").Value = 0
$ws.Range("AM42").Value = 0
$ws.Range("AL43").Value = 0
$ws.Range("AL44").Value = 0
$ws.Range("AM44").Value = 0
$ws.Range("AL45").Value = 0
$ws.Range("AM45").Value = 0
$ws.Range("AL46").Value = 0
$ws.Range("AM46").Value = 0
$ws.Range("AL47").Value = 0
$ws.Range("AM47").Value = 0
$ws.Range("AL48").Value = 0
$ws.Range("AM48").Value = 0
$ws.Range("AL49").Value = 0
$ws.Range("AM49").Value = 0
$ws.Range("AL50").Value = 0
$ws.Range("AM50").Value = 0
$ws.Range("AL51").Value = 0
$ws.Range("AM51").Value = 0
$ws.Range("AL52").Value = 0
$ws.Range("AL53").Value = 0
$ws.Range("AM53").Value = 0
$ws.Range("AI54").Value = 1
$ws.Range("AL54").Value = 0
$ws.Range("AM54").Value = 0
$ws.Range("P54").Value = 1
$ws.Range("T54").Value = 1
$ws.Range("AL55").Value = 0
$ws.Range("AM55").Value = 0
$ws.Range("AL56").Value = 0
$ws.Range("AM56").Value = 0
$ws.Range("AL57").Value = 0
$ws.Range("AM57").Value = 0
$ws.Range("AL58").Value = 0
$ws.Range("AM58").Value = 0
